# cryptos.xlsx update — refresh Price (D) and Volume(1h) (E) columns
# for the coin rows on Sheet1, matching the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.947.14'
$ws.Range('E2').Value = '  +0.14%  '
$ws.Range('D3').Value = '3.118.66'
$ws.Range('E3').Value = '  +2.33%  '
$ws.Range('D4').Value = "'1.00"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = "'560.71"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.72%  '
$ws.Range('D6').Value = "'139.65"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.13%  '
$ws.Range('D7').Value = "'1.00"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.22%  '
$ws.Range('D8').Value = '3.110.45'
$ws.Range('E8').Value = '  +2.34%  '
$ws.Range('D9').Value = "'0.495"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.00%  '
$ws.Range('D10').Value = "'6.78"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +6.33%  '
$ws.Range('D11').Value = "'0.155"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.77%  '
$ws.Range('D12').Value = "'0.456"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.33%  '
$ws.Range('D13').Value = "'35.72"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.68%  '
$ws.Range('D14').Value = "'0.0000218"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.66%  '
$ws.Range('D15').Value = '3.619.87'
$ws.Range('E15').Value = '  +2.29%  '
$ws.Range('D16').Value = '63.938.28'
$ws.Range('E16').Value = '  +0.27%  '
$ws.Range('E17').Value = '  +0.50%  '
$ws.Range('D18').Value = '3.117.25'
$ws.Range('E18').Value = '  +2.10%  '
$ws.Range('D19').Value = "'509.02"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +5.46%  '
$ws.Range('D20').Value = "'6.70"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.42%  '
$ws.Range('D21').Value = "'13.81"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.13%  '
$ws.Range('D22').Value = "'0.710"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +4.66%  '
$ws.Range('D23').Value = "'7.29"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.26%  '
$ws.Range('D24').Value = "'12.45"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.59%  '
$ws.Range('D25').Value = "'78.00"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.68%  '
$ws.Range('D26').Value = "'0.999"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.07%  '
$ws.Range('D27').Value = "'2.79"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +4.13%  '
$ws.Range('D28').Value = "'8.44"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +6.78%  '
$ws.Range('D29').Value = "'2.06"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.06%  '
$ws.Range('D30').Value = "'0.998"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.30%  '
$ws.Range('D31').Value = "'26.34"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.89%  '
$ws.Range('D32').Value = "'2.56"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.67%  '
$ws.Range('E33').Value = '  -0.36%  '
$ws.Range('D34').Value = "'545.28"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -6.19%  '
$ws.Range('D35').Value = "'56.45"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +9.07%  '
$ws.Range('D36').Value = "'5.94"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.16%  '
$ws.Range('D37').Value = "'5.23"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.11%  '
$ws.Range('D38').Value = "'0.0418"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +5.42%  '
$ws.Range('D39').Value = "'0.0803"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.89%  '
$ws.Range('D40').Value = '3.078.25'
$ws.Range('E40').Value = '  +4.76%  '
$ws.Range('D41').Value = "'0.119"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.54%  '
$ws.Range('D42').Value = "'8.15"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.32%  '
$ws.Range('D43').Value = "'2.62"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -7.41%  '
$ws.Range('D44').Value = "'0.257"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +6.78%  '
$ws.Range('E46').Value = '  +1.93%  '
$ws.Range('D47').Value = "'122.01"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.66%  '
$ws.Range('D48').Value = "'24.48"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.29%  '
$ws.Range('E49').Value = '  +0.53%  '
$ws.Range('D50').Value = '0.0₃0504'
$ws.Range('E50').Value = '  -3.35%  '
$ws.Range('E51').Value = '  -0.34%  '
